$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: replace IYR with EEM
$ws.Range("A7").Value = "EEM"

# Row 8: add new GLD row, matching formatting/values of row 7
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)

$ws.Range("A8").Value = "GLD"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.99999000000000005

# Update selection as per diff
$ws.Range("I18").Select()
